# Add a new "Reverse a Linked List" entry to the DW (Data structures/problems)
# worksheet, mirroring the formatting of the row above it ("Lowest common
# ancestor of BST", row 31 - same "Linked List / Easy" style group).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # activeTab = DW sheet

# Row 31 (A31:E31) carries the "Linked List / Easy" look (fonts, fill,
# borders) that the new row should reuse, so copy it down into the new
# row 33 first and then overwrite the copied text with the new content.
$ws.Range("A31:E31").Copy($ws.Range("A33:E33"))

$ws.Range("A33").Value = "Reverse a Linked List"
$ws.Range("B33").Value = "Linked List"
$ws.Range("C33").Value = "Easy"
$ws.Range("D33").Value = "https://leetcode.com/problems/reverse-linked-list/"
$ws.Range("E33").Value = "Use a prev=None pointer and keep changing the next pos"

# Match the row height used by the rest of the table.
$ws.Rows.Item(33).RowHeight = 16

# Leave the freshly-entered row selected (whole row), like the previous
# last row (32) was selected before the edit.
$null = $ws.Rows.Item(33).Select()
